# Apply latest nowcast run: refresh Prognose/Revision figures for existing
# release dates (rows 2-11) and append the new 2025-08-30 row (row 12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"
$ws.Range("A2").Value = "'2025-03-30"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 0.3706856926004305
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("A3").Value = "'2025-04-15"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 0.33760881937391257
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.005324415232551229
$ws.Range("E3").Value = -0.00007063312843090682
$ws.Range("F3").Value = -0.0010370028978019478
$ws.Range("G3").Value = 0.00027158222997034136
$ws.Range("H3").Value = -0.00010133027462725939
$ws.Range("I3").Value = -0.001519019394350117
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.001971711954052935
$ws.Range("A4").Value = "'2025-04-30"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = 0.3248738356501871
$ws.Range("C4").Value = -0.002698130616130056
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.00010830166016220857
$ws.Range("F4").Value = 0.000012842249016542636
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.0002165022596360797
$ws.Range("I4").Value = -0.0011556701473473489
$ws.Range("J4").Value = 0.0001323932788424905
$ws.Range("K4").Value = 0.0010142658407935046
$ws.Range("A5").Value = "'2025-05-15"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 0.33615311833588
$ws.Range("C5").Value = 0.006754961933895529
$ws.Range("D5").Value = -0.007201408935375718
$ws.Range("E5").Value = -0.000032113591087062646
$ws.Range("F5").Value = -0.00008487567670258613
$ws.Range("G5").Value = -0.00035798202084717464
$ws.Range("H5").Value = 0.000004377419650850858
$ws.Range("I5").Value = -0.0007524236511458571
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.00000399043996773818
$ws.Range("A6").Value = "'2025-05-30"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 0.4136507765264203
$ws.Range("C6").Value = 0.032136322593269505
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -0.0002052097839204351
$ws.Range("F6").Value = -0.00004300639169014332
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.00005859509754911743
$ws.Range("I6").Value = -0.0012991789278354172
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -0.0027821348890847974
$ws.Range("A7").Value = "'2025-06-15"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = 0.3776568793456311
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.0021339898911660127
$ws.Range("E7").Value = -0.0005787465538445493
$ws.Range("F7").Value = -0.004608679407095669
$ws.Range("G7").Value = 0.0006488307593035556
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.00015238839535079794
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.0028952513017001524
$ws.Range("A8").Value = "'2025-06-30"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = 0.21427595711212102
$ws.Range("C8").Value = -0.05252765766497771
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = -0.000025056681412530756
$ws.Range("F8").Value = -0.0004571646341500491
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.00005356520276720551
$ws.Range("I8").Value = 0.002079745182253993
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = -0.0001455471442038947
$ws.Range("A9").Value = "'2025-07-15"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = 0.2467482416078785
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.004256175895418775
$ws.Range("E9").Value = -0.00366990378984903
$ws.Range("F9").Value = -0.005100119231592561
$ws.Range("G9").Value = 0.0008005756618776729
$ws.Range("H9").Value = -0.0001500153511616265
$ws.Range("I9").Value = -0.0005680210641357379
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0.0006956746699531557
$ws.Range("A10").Value = "'2025-07-30"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = 0.43607381961049624
$ws.Range("C10").Value = 0.08785766720347168
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -0.00022451145117720583
$ws.Range("F10").Value = -0.0004473778018075806
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = -0.000020509517765437424
$ws.Range("I10").Value = -0.0011270707570654214
$ws.Range("J10").Value = -0.002319187173705693
$ws.Range("K10").Value = 0.0020990652586068415
$ws.Range("A11").Value = "'2025-08-15"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = 0.3557322547826534
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = -0.022222694135920793
$ws.Range("E11").Value = 0.002103000119235687
$ws.Range("F11").Value = 0.0035123932202059
$ws.Range("G11").Value = 0.0014918558862453454
$ws.Range("H11").Value = 0.0003419645011809428
$ws.Range("I11").Value = -0.003877602971283626
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0.007077918722477927
$ws.Range("A12").Value = "'2025-08-30"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = 0.2888413589105628
$ws.Range("C12").Value = -0.052272272985380425
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.0001522051420676436
$ws.Range("F12").Value = 0.00003102292354173377
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0.00002236659085877483
$ws.Range("I12").Value = 0.0006144139477576444
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = -0.00720223628004607
